$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 448845.53
$ws.Range("I17").Value = 99
$ws.Range("J17").Value = 506748.3
$ws.Range("K17").Value = 297
$ws.Range("L17").Value = 1520244.9
$ws.Range("M17").Value = -129
$ws.Range("N17").Value = -1520580.9

$ws.Range("H74").Value = 4033.2917
$ws.Range("I74").Value = 2973.2222
$ws.Range("J74").Value = 4669.3335
$ws.Range("K74").Value = 2973.2222
$ws.Range("L74").Value = 4669.3335
$ws.Range("M74").Value = -2037.2222
$ws.Range("N74").Value = -6541.3335

$ws.Range("H77").Value = 4033.2917
$ws.Range("I77").Value = 2973.2222
$ws.Range("J77").Value = 4669.3335
$ws.Range("K77").Value = 14866.111
$ws.Range("L77").Value = 23346.6675
$ws.Range("M77").Value = -10186.111
$ws.Range("N77").Value = -32706.6675

$ws.Range("H131").Value = 2893.4783
$ws.Range("I131").Value = 1560.3572
$ws.Range("J131").Value = 4967.222
$ws.Range("K131").Value = 4681.071599999999
$ws.Range("L131").Value = 14901.666
$ws.Range("M131").Value = 358.9284000000007
$ws.Range("N131").Value = -24981.666

$ws.Range("H132").Value = 3447.2654
$ws.Range("I132").Value = 2847.6667
$ws.Range("J132").Value = 5107.6924
$ws.Range("K132").Value = 8543.000100000001
$ws.Range("L132").Value = 15323.0772
$ws.Range("M132").Value = -6013.000100000001
$ws.Range("N132").Value = -20383.0772

$ws.Range("H138").Value = 2458.3167
$ws.Range("I138").Value = 1098.5483
$ws.Range("J138").Value = 3911.862
$ws.Range("K138").Value = 3295.6449
$ws.Range("L138").Value = 11735.586
$ws.Range("M138").Value = 1844.3551
$ws.Range("N138").Value = -22015.586

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 325.625
$ws.Range("I5").Value = 170.5
$ws.Range("J5").Value = 377.33334
$ws.Range("K5").Value = 170.5
$ws.Range("L5").Value = 377.33334
$ws.Range("M5").Value = -58.5
$ws.Range("N5").Value = -601.33334

$ws.Range("H61").Value = 1731.4
$ws.Range("I61").Value = 1463.1904
$ws.Range("J61").Value = 2357.2222
$ws.Range("K61").Value = 1463.1904
$ws.Range("L61").Value = 2357.2222
$ws.Range("M61").Value = -1251.1904
$ws.Range("N61").Value = -2781.2222

$ws.Range("H63").Value = 2153
$ws.Range("I63").Value = 0
$ws.Range("J63").Value = 2153
$ws.Range("K63").Value = 0
$ws.Range("L63").Value = 2153
$ws.Range("N63").Value = -3525
$ws.Range("M63").ClearContents()

$ws.Range("H66").Value = 2153
$ws.Range("I66").Value = 0
$ws.Range("J66").Value = 2153
$ws.Range("K66").Value = 0
$ws.Range("L66").Value = 10765
$ws.Range("N66").Value = -17629
$ws.Range("M66").ClearContents()

$ws.Range("H136").Value = 1731.4
$ws.Range("I136").Value = 1463.1904
$ws.Range("J136").Value = 2357.2222
$ws.Range("K136").Value = 4389.5712
$ws.Range("L136").Value = 7071.6666
$ws.Range("M136").Value = -1839.5712
$ws.Range("N136").Value = -12171.6666

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 325.625
$ws.Range("I4").Value = 170.5
$ws.Range("J4").Value = 377.33334
$ws.Range("K4").Value = 170.5
$ws.Range("L4").Value = 377.33334
$ws.Range("M4").Value = -55.5
$ws.Range("N4").Value = -607.33334

$ws.Range("H82").Value = 6714.6
$ws.Range("I82").Value = 822.5
$ws.Range("J82").Value = 30283
$ws.Range("K82").Value = 822.5
$ws.Range("L82").Value = 30283
$ws.Range("M82").Value = -439.5
$ws.Range("N82").Value = -31049

$ws.Range("H85").Value = 6714.6
$ws.Range("I85").Value = 822.5
$ws.Range("J85").Value = 30283
$ws.Range("K85").Value = 822.5
$ws.Range("L85").Value = 30283
$ws.Range("M85").Value = 503.5
$ws.Range("N85").Value = -32935

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 40085.57
$ws.Range("I31").Value = 1369.7241
$ws.Range("J31").Value = 81669.25999999999
$ws.Range("K31").Value = 1369.7241
$ws.Range("L31").Value = 81669.25999999999
$ws.Range("M31").Value = -1074.7241
$ws.Range("N31").Value = -82259.25999999999

$ws.Range("H34").Value = 40085.57
$ws.Range("I34").Value = 1369.7241
$ws.Range("J34").Value = 81669.25999999999
$ws.Range("K34").Value = 1369.7241
$ws.Range("L34").Value = 81669.25999999999
$ws.Range("M34").Value = -1167.7241
$ws.Range("N34").Value = -82073.25999999999

$ws.Range("H107").Value = 48509.523
$ws.Range("I107").Value = 84105.914
$ws.Range("J107").Value = 1047.6666
$ws.Range("K107").Value = 84105.914
$ws.Range("L107").Value = 1047.6666
$ws.Range("M107").Value = -82185.914
$ws.Range("N107").Value = -4887.6666

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 897.0351000000001
$ws.Range("I131").Value = 609.1579
$ws.Range("J131").Value = 1040.9736
$ws.Range("K131").Value = 1827.4737
$ws.Range("L131").Value = 3122.9208
$ws.Range("M131").Value = 3212.5263
$ws.Range("N131").Value = -13202.9208

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 6471.722
$ws.Range("I102").Value = 7274.909
$ws.Range("J102").Value = 5209.5713
$ws.Range("K102").Value = 7274.909
$ws.Range("L102").Value = 5209.5713
$ws.Range("M102").Value = -5652.909
$ws.Range("N102").Value = -8453.5713

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 1510.742
$ws.Range("I61").Value = 1182.1538
$ws.Range("J61").Value = 1748.0555
$ws.Range("K61").Value = 1182.1538
$ws.Range("L61").Value = 1748.0555
$ws.Range("M61").Value = -980.1538
$ws.Range("N61").Value = -2152.0555

$ws.Range("H74").Value = 39500
$ws.Range("I74").Value = 0
$ws.Range("J74").Value = 39500
$ws.Range("K74").Value = 0
$ws.Range("L74").Value = 39500
$ws.Range("N74").Value = -41496

$ws.Range("H77").Value = 39500
$ws.Range("I77").Value = 0
$ws.Range("J77").Value = 39500
$ws.Range("K77").Value = 0
$ws.Range("L77").Value = 118500
$ws.Range("N77").Value = -128484

$ws.Range("H113").Value = 1510.742
$ws.Range("I113").Value = 1182.1538
$ws.Range("J113").Value = 1748.0555
$ws.Range("K113").Value = 1182.1538
$ws.Range("L113").Value = 1748.0555
$ws.Range("M113").Value = 987.8462
$ws.Range("N113").Value = -6088.0555

$ws.Range("H122").Value = 41053.117
$ws.Range("I122").Value = 49865.76
$ws.Range("J122").Value = 4040
$ws.Range("K122").Value = 149597.28
$ws.Range("L122").Value = 12120
$ws.Range("M122").Value = -147147.28
$ws.Range("N122").Value = -17020

$ws.Range("H132").Value = 6016.5713
$ws.Range("I132").Value = 10111.733
$ws.Range("J132").Value = 2945.2
$ws.Range("K132").Value = 30335.199
$ws.Range("L132").Value = 8835.599999999999
$ws.Range("M132").Value = -27805.199
$ws.Range("N132").Value = -13895.6

$ws.Range("H136").Value = 3117.5789
$ws.Range("I136").Value = 1285.0454
$ws.Range("J136").Value = 9320
$ws.Range("K136").Value = 3855.1362
$ws.Range("L136").Value = 27960
$ws.Range("M136").Value = -1305.1362
